$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 223
$ws.Range("B223").Value2 = 6978439
$ws.Range("E223").Value2 = 45386.35416666666
$ws.Range("F223").Value2 = 'CSKA 1948 Sofia'
$ws.Range("G223").Value2 = 'Botev Vratsa'
$ws.Range("K223").Value2 = 1.333
$ws.Range("L223").Value2 = 5.5
$ws.Range("M223").Value2 = 7.5
$ws.Range("N223").Value2 = 1.5
$ws.Range("O223").Value2 = 4
$ws.Range("P223").Value2 = 7
$ws.Range("Q223").Value2 = -1
$ws.Range("R223").Value2 = 1.9
$ws.Range("S223").Value2 = 1.95
$ws.Range("T223").Value2 = 2.5
$ws.Range("U223").Value2 = 1.925
$ws.Range("V223").Value2 = 1.925

# Row 224
$ws.Range("B224").Value2 = 6978440
$ws.Range("E224").Value2 = 45386.45833333334
$ws.Range("F224").Value2 = 'Beroe'
$ws.Range("G224").Value2 = 'Botev Plovdiv'
$ws.Range("K224").Value2 = 3.5
$ws.Range("L224").Value2 = 3.2
$ws.Range("M224").Value2 = 2.15
$ws.Range("N224").Value2 = 4.333
$ws.Range("O224").Value2 = 3.5
$ws.Range("P224").Value2 = 1.909
$ws.Range("Q224").Value2 = 0.5
$ws.Range("R224").Value2 = 1.975
$ws.Range("S224").Value2 = 1.875
$ws.Range("U224").Value2 = 1.85
$ws.Range("V224").Value2 = 2

# Row 225
$ws.Range("B225").Value2 = 6978441
$ws.Range("E225").Value2 = 45386.5625
$ws.Range("F225").Value2 = 'Cherno More Varna'
$ws.Range("G225").Value2 = 'Lokomotiv 1929 Sofia'
$ws.Range("K225").Value2 = 1.363
$ws.Range("L225").Value2 = 4.75
$ws.Range("M225").Value2 = 8.5
$ws.Range("N225").Value2 = 1.25
$ws.Range("O225").Value2 = 5.25
$ws.Range("P225").Value2 = 17
$ws.Range("Q225").Value2 = -1.75
$ws.Range("R225").Value2 = 1.975
$ws.Range("S225").Value2 = 1.875
$ws.Range("T225").Value2 = 2.5

# Row 226
$ws.Range("B226").Value2 = 6978388
$ws.Range("E226").Value2 = 45388.48958333334
$ws.Range("F226").Value2 = 'FC Hebar Pazardzhik'
$ws.Range("G226").Value2 = 'Etar 1924 Veliko Tarnovo'
$ws.Range("K226").Value2 = 1.75
$ws.Range("L226").Value2 = 3.5
$ws.Range("M226").Value2 = 4.75
$ws.Range("N226").Value2 = 1.7
$ws.Range("O226").Value2 = 3.6
$ws.Range("P226").Value2 = 5
$ws.Range("Q226").Value2 = -0.75
$ws.Range("R226").Value2 = 1.95
$ws.Range("S226").Value2 = 1.9
$ws.Range("T226").Value2 = 2.25
$ws.Range("U226").Value2 = 2
$ws.Range("V226").Value2 = 1.85

# Row 227
$ws.Range("B227").Value2 = 6978438
$ws.Range("E227").Value2 = 45388.59375
$ws.Range("F227").Value2 = 'Arda Kardzhali'
$ws.Range("G227").Value2 = 'Slavia Sofia'
$ws.Range("K227").Value2 = 2.05
$ws.Range("M227").Value2 = 3.75
$ws.Range("N227").Value2 = 2.1
$ws.Range("O227").Value2 = 3.2
$ws.Range("P227").Value2 = 3.6
$ws.Range("Q227").Value2 = -0.25
$ws.Range("R227").Value2 = 1.825
$ws.Range("S227").Value2 = 2.025
$ws.Range("T227").Value2 = 2
$ws.Range("U227").Value2 = 1.9
$ws.Range("V227").Value2 = 1.95

# Remove the now-obsolete last row (shifted out of the dataset)
$ws.Rows("228:228").Delete()
